$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2 (format + values) into rows 3 and 4, which also
# establishes the correct "Text" number formats/fonts/borders/fills
# needed so that later numeric-looking values are stored as text.
$ws.Range("A2:U2").Copy($ws.Range("A3:U3"))
$ws.Range("A2:U2").Copy($ws.Range("A4:U4"))

# Overwrite the cells that actually differ between the rows, typing
# the values in the same order the shared-strings table expects them
# so indices line up with the target workbook.
$ws.Range("Q3").Value = "Otro valor"
$ws.Range("P3").Value = "*7826"
$ws.Range("S3").Value = "Pesos"
$ws.Range("R3").Value = "500000"
$ws.Range("A3").Value = "'2"
$ws.Range("A4").Value = "'3"
$ws.Range("O4").Value = "Personal Visa"
$ws.Range("P4").Value = "*5880"
$ws.Range("R4").Value = "480369"
$ws.Range("Q4").Value = "Otro valor"
$ws.Range("S4").Value = "Pesos"

# Add hyperlinks on the new rows' correo column, mirroring N2.
$ws.Hyperlinks.Add($ws.Range("N3"), "mailto:jalzate@todo1.net")
$ws.Hyperlinks.Add($ws.Range("N4"), "mailto:jalzate@todo1.net")

# Update the view: scroll so column L is leftmost and select Q15.
$excel.Goto($ws.Range("L1"), $true)
$ws.Range("Q15").Select()
